$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (column D) and "Volume(1h)" (column E) figures for the
# crypto list with the latest scraped values. All cells in this sheet hold
# plain text (e.g. "60.828.24", "  -0.70%  "), so for any new Price value
# that looks like a genuine number (e.g. "404.21") we first force the cell's
# NumberFormat to Text ("@") - otherwise Excel's COM layer would silently
# reinterpret the assignment as a numeric value and drop the original text
# representation (trailing zeros, thousand-dot grouping, etc.).
$ws.Range("D2").Value = "61.138.68"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.376.51"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "404.21"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.34"
$ws.Range("E6").Value = "  +13.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("E7").Value = "  +7.15%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.673"
$ws.Range("E9").Value = "  +8.82%  "
$ws.Range("E10").Value = "  +11.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.11"
$ws.Range("E11").Value = "  +8.96%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "3.919.73"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.50"
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.62"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "3.350.03"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.40"
$ws.Range("E17").Value = "  +8.60%  "
$ws.Range("D18").Value = "60.927.64"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000134"
$ws.Range("E20").Value = "  +19.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.24"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "82.62"
$ws.Range("E22").Value = "  +13.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.16"
$ws.Range("E23").Value = "  +7.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "306.69"
$ws.Range("E24").Value = "  +4.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.73"
$ws.Range("E26").Value = "  +5.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.46"
$ws.Range("E27").Value = "  +14.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.59"
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.44"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").Value = "  +6.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.79"
$ws.Range("E32").Value = "  +6.43%  "
$ws.Range("E33").Value = "  +6.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.53"
$ws.Range("E34").Value = "  +7.31%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0485"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.31"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.995"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.03"
$ws.Range("E41").Value = "  +8.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.126"
$ws.Range("E42").Value = "  +5.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.02"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.95"
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.284"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.95"
$ws.Range("E46").Value = "  +5.21%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.82"
$ws.Range("E48").Value = "  +4.78%  "
$ws.Range("D49").Value = "2.137.09"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").Value = "3.702.63"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  +1.48%  "
